$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the commit diff (price refresh + row43/41 coin swap)
$updates = @{
    'D2' = '246.23'
    'G2' = '8'
    'G3' = '8'
    'D4' = '5.461'
    'G4' = '8'
    'D5' = '0.05672'
    'G5' = '8'
    'D6' = '3.373'
    'G6' = '8'
    'D7' = '0.8013'
    'G7' = '8'
    'G8' = '8'
    'D9' = '0.1428'
    'G9' = '8'
    'D10' = '0.07311'
    'G10' = '8'
    'D11' = '0.03165'
    'G11' = '8'
    'D12' = '0.02943'
    'G12' = '8'
    'D13' = '0.09285'
    'G13' = '8'
    'D14' = '0.001651'
    'G14' = '8'
    'D15' = '3.218'
    'G15' = '8'
    'D16' = '0.04699'
    'G16' = '8'
    'D17' = '0.0005892'
    'E17' = '16OneONEWorstin24h'
    'G17' = '8'
    'D18' = '0.006345'
    'G18' = '8'
    'D19' = '0.005041'
    'G19' = '8'
    'D20' = '0.001047'
    'G20' = '8'
    'D21' = '0.0001502'
    'G21' = '8'
    'D22' = '0.0003201'
    'G22' = '8'
    'D23' = '3.805'
    'G23' = '8'
    'D24' = '6.423'
    'G24' = '8'
    'D25' = '2.112'
    'G25' = '8'
    'D26' = '0.3328'
    'G26' = '8'
    'G27' = '8'
    'G28' = '8'
    'G29' = '8'
    'G30' = '8'
    'G31' = '8'
    'G32' = '8'
    'G33' = '8'
    'G34' = '8'
    'G35' = '8'
    'G36' = '8'
    'G37' = '8'
    'G38' = '8'
    'G39' = '8'
    'D40' = '0.04082'
    'G40' = '8'
    'B41' = 'KickToken'
    'C41' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D41' = '0.006957'
    'E41' = '40KickTokenKICK'
    'G41' = '8'
    'D42' = '0.003504'
    'G42' = '8'
    'B43' = 'BKEXToken'
    'C43' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D43' = '0.1040'
    'E43' = '42BKEXTokenBKK'
    'G43' = '8'
    'D44' = '0.008093'
    'G44' = '8'
    'D45' = '0.00005849'
    'G45' = '8'
    'G46' = '8'
    'D47' = '0.6827'
    'G47' = '8'
    'D48' = '0.01023'
    'E48' = '47BOLOBOLO'
    'G48' = '8'
    'D49' = '0.00002102'
    'G49' = '8'
    'G50' = '8'
    'G51' = '8'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
